# Update the "dSF" (column F) values on Sheet1 to reflect the repulled data.
# These changes only touch column F; column E ("dS0") is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = 2
$ws.Range("F5").Value = 2
$ws.Range("F12").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("F22").Value = -5
$ws.Range("F23").Value = 1
$ws.Range("F25").Value = -3
$ws.Range("F28").Value = -4
$ws.Range("F32").Value = -4
$ws.Range("F37").Value = -12
$ws.Range("F38").Value = -1
$ws.Range("F41").Value = -3
$ws.Range("F43").Value = -7
$ws.Range("F46").Value = 2
$ws.Range("F47").Value = 2
$ws.Range("F52").Value = -6
$ws.Range("F53").Value = -1
$ws.Range("F59").Value = -3
$ws.Range("F62").Value = 4
$ws.Range("F66").Value = -1
$ws.Range("F69").Value = -5
